$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: update title and link
$ws.Range("D4").Value = "#02-파이썬(Python) 리스트(list)와 튜플(tuple)"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-02"

# Row 42: update title and link
$ws.Range("D42").Value = "GetPrivateProfile 실패"
$ws.Range("E42").Value = "https://kjk92.tistory.com/69"

# Row 51: update title and link
$ws.Range("D51").Value = "[javascript] 문자열을 정수형 또는 실수형으로 변환하려면, parseInt(), parseFloat()"
$ws.Range("E51").Value = "https://bskyvision.com/1198"
